$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update simulated transition-matrix probabilities (more games simulated)
$ws.Range("B2").Value = 0.1965217391304348
$ws.Range("C2").Value = 0.5478260869565217
$ws.Range("J2").Value = 0.01217391304347826
$ws.Range("P2").Value = 0.1686956521739131
$ws.Range("S2").Value = 0.07478260869565218
$ws.Range("B3").Value = 0.008746355685131196
$ws.Range("C3").Value = 0.03206997084548105
$ws.Range("J3").Value = 0.03790087463556852
$ws.Range("P3").Value = 0.7376093294460642
$ws.Range("S3").Value = 0.1836734693877551
$ws.Range("J4").Value = 0.09183673469387756
$ws.Range("O4").Value = 0.02040816326530612
$ws.Range("P4").Value = 0.5
$ws.Range("S4").Value = 0.3877551020408163
$ws.Range("J5").Value = 0.3333333333333333
$ws.Range("P5").Value = 0.3333333333333333
$ws.Range("S5").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.0611764705882353
$ws.Range("D6").Value = 0.01411764705882353
$ws.Range("F6").Value = 0.05411764705882353
$ws.Range("J6").Value = 0.2235294117647059
$ws.Range("O6").Value = 0.02117647058823529
$ws.Range("Q6").Value = 0.1835294117647059
$ws.Range("R6").Value = 0.06588235294117648
$ws.Range("B7").Value = 0.1155015197568389
$ws.Range("D7").Value = 0.0243161094224924
$ws.Range("E7").Value = 0.00303951367781155
$ws.Range("F7").Value = 0.06382978723404255
$ws.Range("J7").Value = 0.121580547112462
$ws.Range("O7").Value = 0.00911854103343465
$ws.Range("Q7").Value = 0.1945288753799392
$ws.Range("R7").Value = 0.05167173252279635
$ws.Range("S7").Value = 0.4164133738601823
$ws.Range("B8").Value = 0.08009153318077804
$ws.Range("D8").Value = 0.01830663615560641
$ws.Range("E8").Value = 0.0011441647597254
$ws.Range("F8").Value = 0.06636155606407322
$ws.Range("J8").Value = 0.1018306636155606
$ws.Range("O8").Value = 0.01487414187643021
$ws.Range("Q8").Value = 0.2162471395881007
$ws.Range("R8").Value = 0.07894736842105263
$ws.Range("S8").Value = 0.4221967963386727
$ws.Range("B9").Value = 0.1007556675062972
$ws.Range("D9").Value = 0.03526448362720403
$ws.Range("F9").Value = 0.07052896725440806
$ws.Range("J9").Value = 0.08312342569269521
$ws.Range("O9").Value = 0.007556675062972292
$ws.Range("Q9").Value = 0.2241813602015113
$ws.Range("R9").Value = 0.07556675062972293
$ws.Range("S9").Value = 0.4030226700251889
$ws.Range("B10").Value = 0.1232394366197183
$ws.Range("D10").Value = 0.02332746478873239
$ws.Range("E10").Value = 0.0008802816901408451
$ws.Range("F10").Value = 0.07614436619718309
$ws.Range("J10").Value = 0.1025528169014085
$ws.Range("O10").Value = 0.02024647887323944
$ws.Range("Q10").Value = 0.2143485915492958
$ws.Range("R10").Value = 0.07614436619718309
$ws.Range("S10").Value = 0.3631161971830986
$ws.Range("G11").Value = 0.0990990990990991
$ws.Range("J11").Value = 0.1081081081081081
$ws.Range("K11").Value = 0.1351351351351351
$ws.Range("L11").Value = 0.6418918918918919
$ws.Range("S11").Value = 0.01576576576576576
$ws.Range("G12").Value = 0.7892976588628763
$ws.Range("J12").Value = 0.1538461538461539
$ws.Range("K12").Value = 0.006688963210702341
$ws.Range("L12").Value = 0.03678929765886288
$ws.Range("S12").Value = 0.01337792642140468
$ws.Range("G13").Value = 0.6086956521739131
$ws.Range("J13").Value = 0.3260869565217391
$ws.Range("S13").Value = 0.06521739130434782
$ws.Range("G14").Value = 0.8
$ws.Range("J14").Value = 0.2
$ws.Range("F15").Value = 0.0137299771167048
$ws.Range("H15").Value = 0.1876430205949657
$ws.Range("I15").Value = 0.06864988558352403
$ws.Range("J15").Value = 0.3684210526315789
$ws.Range("K15").Value = 0.05034324942791762
$ws.Range("M15").Value = 0.011441647597254
$ws.Range("O15").Value = 0.06636155606407322
$ws.Range("S15").Value = 0.2334096109839817
$ws.Range("F16").Value = 0.02406417112299465
$ws.Range("H16").Value = 0.1550802139037433
$ws.Range("I16").Value = 0.09090909090909091
$ws.Range("J16").Value = 0.4117647058823529
$ws.Range("K16").Value = 0.1016042780748663
$ws.Range("M16").Value = 0.0374331550802139
$ws.Range("O16").Value = 0.04545454545454546
$ws.Range("S16").Value = 0.1336898395721925
$ws.Range("F17").Value = 0.01416122004357298
$ws.Range("H17").Value = 0.1840958605664488
$ws.Range("I17").Value = 0.1045751633986928
$ws.Range("J17").Value = 0.420479302832244
$ws.Range("K17").Value = 0.07516339869281045
$ws.Range("M17").Value = 0.02287581699346405
$ws.Range("N17").Value = 0.001089324618736384
$ws.Range("O17").Value = 0.08169934640522876
$ws.Range("S17").Value = 0.09586056644880174
$ws.Range("F18").Value = 0.02515723270440252
$ws.Range("H18").Value = 0.1886792452830189
$ws.Range("I18").Value = 0.1226415094339623
$ws.Range("J18").Value = 0.4025157232704403
$ws.Range("K18").Value = 0.07861635220125786
$ws.Range("M18").Value = 0.01572327044025157
$ws.Range("N18").Value = 0.003144654088050315
$ws.Range("O18").Value = 0.09119496855345911
$ws.Range("S18").Value = 0.07232704402515723
$ws.Range("F19").Value = 0.01424121050289275
$ws.Range("H19").Value = 0.225634178905207
$ws.Range("I19").Value = 0.09479305740987984
$ws.Range("J19").Value = 0.3751668891855808
$ws.Range("K19").Value = 0.1014686248331108
$ws.Range("M19").Value = 0.02225189141076991
$ws.Range("N19").Value = 0.001335113484646195
$ws.Range("O19").Value = 0.07610146862483311
$ws.Range("S19").Value = 0.09479305740987984

Write-Output "Updated 117 cells"
